# Sprint 2 Hours Log -- add two more logged entries (Colors / SF-15) and
# move the selection down to D4, matching the author's "Updated my hours
# log" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Colors -----------------------------------------------------
$ws.Range("A3").Value = Get-Date -Year 2017 -Month 4 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Colors"
$ws.Range("D3").Value = "Added colors just for the fun of it"

# --- Row 4: SF-15 --------------------------------------------------------
$ws.Range("A4").Value = Get-Date -Year 2017 -Month 4 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Range("B4").Value = 0.5
# Write the Description before the User Story text so the shared-string
# table fills in the same order the workbook author produced.
$ws.Range("D4").Value = "Added the menu, menu item, and empty action listener for the Top 10 List"
$ws.Range("C4").Value = "SF-15"

# Match the look of the existing logged row (row 2): centered, bordered,
# 12pt Calibri cells for the User Story / Description columns.
$ws.Range("C2:D2").Copy()
$ws.Range("C3:D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection on D4, same as the saved workbook.
$ws.Range("D4").Select() | Out-Null
